$d = $word.ActiveDocument

# Locate the edit point: right after "...die entnommen werden soll" and
# before " jedoch nicht mehr vorhanden ist, wird eine Warnung ausgegeben."
$anchor = "werden soll"
$searchRange = $d.Content
$searchRange.Find.Execute("$anchor jedoch nicht mehr vorhanden ist")
if ($searchRange.Find.Found) {
    $insertPos = $searchRange.Start + $anchor.Length

    # Word keeps the "_GoBack" bookmark pinned to the location of the most
    # recent edit. Drop the old one (currently sitting at the top of the
    # document) before we make the new edit.
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks.Item("_GoBack").Delete()
    }

    # Type the missing comma at the located position.
    $insRange = $d.Range($insertPos, $insertPos)
    $insRange.InsertAfter(",")

    # Re-create "_GoBack" collapsed right after the newly typed comma,
    # matching Word's behaviour of tracking the last edit location.
    $newGoBackRange = $d.Range($insertPos + 1, $insertPos + 1)
    $d.Bookmarks.Add("_GoBack", $newGoBackRange)
}
